$wb = $excel.ActiveWorkbook

# Use the "Germany" sheet as the template for the new country sheets - it
# has the exact column widths / styles the new country sheets end up with.
$template = $wb.Worksheets.Item("Germany")

# --- Netherlands ------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$template.Copy($null, $lastSheet)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2199"
$netherlands.Range("B2").Value = "Netherlands Market"
[void]$netherlands.Activate()
[void]$netherlands.Range("B4").Select()

# --- Austria ------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$template.Copy($null, $lastSheet)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/T2306"
$austria.Range("B2").Value = "Austria Market"
[void]$austria.Activate()
[void]$austria.Range("C14").Select()

# --- Denmark --------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$template.Copy($null, $lastSheet)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2798"
$denmark.Range("B2").Value = "Denmark Market"
[void]$denmark.Activate()
[void]$denmark.Range("D22").Select()

# Austria's "Wg" row (A9) uses a different constant than the template -
# this was edited after all three country sheets already existed.
$austria.Range("A9").Value = "Fire Brigade Panel"

# Austria ends up being the active sheet / selected tab.
[void]$austria.Activate()
[void]$austria.Range("C14").Select()
